$wb = $excel.ActiveWorkbook

# ---- Sheet "Resumen": update maximum time value ----
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("C2").Value = 505.47450535193

# ---- Sheet "Solucion": reassign Pedido/Salida pairs (new nearest-neighbor route) ----
$wsSolucion = $wb.Worksheets.Item("Solucion")

$pedidos = @(
    "Pedido_16","Pedido_28","Pedido_35","Pedido_5","Pedido_7","Pedido_15","Pedido_22","Pedido_33","Pedido_6","Pedido_38",
    "Pedido_2","Pedido_17","Pedido_3","Pedido_25","Pedido_40","Pedido_14","Pedido_18","Pedido_37","Pedido_19","Pedido_4",
    "Pedido_13","Pedido_9","Pedido_23","Pedido_11","Pedido_39","Pedido_12","Pedido_10","Pedido_27","Pedido_36","Pedido_26",
    "Pedido_32","Pedido_1","Pedido_20","Pedido_30","Pedido_24","Pedido_31","Pedido_21","Pedido_29","Pedido_34","Pedido_8"
)

$salidas = @(
    "S001","S021","S031","S011","S022","S002","S032","S012","S023","S003",
    "S013","S033","S024","S004","S034","S014","S025","S005","S035","S015",
    "S026","S006","S016","S036","S007","S027","S017","S037","S008","S028",
    "S018","S038","S009","S029","S019","S039","S010","S030","S040","S020"
)

for ($i = 0; $i -lt $pedidos.Length; $i++) {
    $row = $i + 2
    $wsSolucion.Cells.Item($row, 1).Value = $pedidos[$i]
    $wsSolucion.Cells.Item($row, 2).Value = $salidas[$i]
}

# ---- Sheet "Metricas": update per-zone times ----
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 505.47450535193
$wsMetricas.Range("B3").Value = 503.6471996972646
